# "Asignar Empleado" use-case doc: in the "Flujo alternativo" row, the
# second cell used to read a numbered list item "Dar de alta empleado".
# It must become a plain paragraph with just "N/A" (no list style / numbering),
# while keeping the existing line-spacing formatting.

$d = $word.ActiveDocument

# Locate the "Flujo alternativo" row robustly via the Tables collection
# (first/only table in the doc) rather than a hard-coded paragraph index.
$table = $d.Tables.Item(1)
$targetRow = 0
for ($r = 1; $r -le $table.Rows.Count; $r++) {
    # Cell text carries trailing cell-mark/paragraph-mark control chars
    # (CR + cell-end marker), so use StartsWith rather than exact equality.
    $label = $table.Cell($r, 1).Range.Text
    if ($label.StartsWith("Flujo alternativo")) {
        $targetRow = $r
        break
    }
}

$cell = $table.Cell($targetRow, 2)
$para = $cell.Range.Paragraphs.Item(1)

# Replace the paragraph's contents/formatting with the exact target markup:
# a plain paragraph (no pStyle, no numPr) that keeps the spacing override
# and contains a single run reading "N/A". InsertXML replaces the range's
# contents in place, preserving the paragraph's own identity attributes.
$xml = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" ' +
       'xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml" ' +
       'w14:paraId="1656CBD3" w14:textId="3D8F598D" w:rsidR="00EE20B7" ' +
       'w:rsidRDefault="00E77BE0" w:rsidP="00EE20B7">' +
       '<w:pPr><w:spacing w:line="240" w:lineRule="auto"/></w:pPr>' +
       '<w:r><w:t>N/A</w:t></w:r>' +
       '</w:p>'

$null = $para.Range.InsertXML($xml)
